# Adds a "localization" column (E) populated with "cyt" for every data row,
# adds a trailing "#Measurement" section header row, and cleans up the
# now-unused fourth custom style that used to be applied to the last
# reaction block (rows 49-60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header ---
$ws.Range("E1").Value = "localization"

# --- Data rows that get "cyt" with NO extra styling (cells stay plain / style 0) ---
$plainRows = @(3,5,7,9,11,13,15,17,19,21,23,25,27,29,31,33,34,36,40,42)
foreach ($r in $plainRows) {
    $ws.Range("E$r").Value = "cyt"
}

# --- Data rows that get "cyt" styled like the sheet's alternate-row style (cellXf 1) ---
$styledRows = @(4,6,8,10,14,16,18,20,22,26,28,30,37,39,41,43,45,49,50,51,52,53,54,55,56,57,58,59,60)
foreach ($r in $styledRows) {
    $c = $ws.Range("E$r")
    $c.Value = "cyt"
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
}

# --- Row 32 is a section-header row; it only picks up the (empty) styled cell in col E ---
$e32 = $ws.Range("E32")
$e32.Font.Name = "Calibri"
$e32.Font.Size = 11

# --- New trailing section header row ---
$ws.Range("A61").Value = "#Measurement"

# --- The old style previously on A49:D60 (cellXf index 4) is no longer used;
#     clear it back to the default/general style so those cells drop their "s" attribute. ---
$ws.Range("A49:D60").ClearFormats()

# --- Match the author's final selection/viewport state ---
$ws.Range("E34").Select()
